$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A3 text (ratingCS-1 -> ratingCS+3)
$ws.Range("A3").Value = "Instructions_EN/ratingCS+3.png"

# Add new row 4 (ratingCS+4)
$ws.Range("A4").Value = "Instructions_EN/ratingCS+4.png"
$ws.Range("B4").Value = -0.1

# Copy formatting of row 3 (A3:B3) down into row 4 (A4:B4) so it starts from the
# same base font styling before the new shading is applied.
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A4").RowHeight = 16

# Apply new light-gray fill banding to rows 2 and 4
$grayRange = $ws.Range("A2:B2,A4:B4")
$grayRange.Interior.ThemeColor = 1
$grayRange.Interior.TintAndShade = -0.049989318521683403

$ws.Range("A2").RowHeight = 15
